# Applies the cryptos-list price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '69.237.90'
$ws.Range('E2').Value = '  +2.05%  '
$ws.Range('D3').Value = '3.317.85'
$ws.Range('E3').Value = '  +2.23%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '591.68'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.18%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '187.85'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.53%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.607'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +1.78%  '
$ws.Range('E9').Value = '  +5.43%  '
$ws.Range('E10').Value = '  -0.88%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.425'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.73%  '
$ws.Range('D12').Value = '3.894.90'
$ws.Range('E12').Value = '  +2.44%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.137'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.20%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '29.38'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +4.83%  '
$ws.Range('D15').Value = '69.243.27'
$ws.Range('E15').Value = '  +2.30%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000174'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +3.81%  '
$ws.Range('D17').Value = '3.334.85'
$ws.Range('E17').Value = '  +2.86%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.93'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.60%  '
$ws.Range('E19').Value = '  +3.28%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '390.29'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +4.09%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.85'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +3.31%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '72.11'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.09%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('E24').Value = '  +3.69%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.522'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.57%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.97'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.70%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.191'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +5.65%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.998'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.50%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.92'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +5.14%  '
$ws.Range('E30').Value = '  +1.79%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '23.22'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.59%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.33'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +3.80%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.28'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +5.79%  '
$ws.Range('E35').Value = '  +4.93%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '163.56'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.95%  '
$ws.Range('E37').Value = '  +2.72%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.844'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.45%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '27.02'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.24%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.68'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +3.04%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.65'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +5.48%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.67'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.55%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '25.93'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.45%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0701'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +3.87%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '41.66'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +2.83%  '
$ws.Range('D46').Value = '2.668.21'
$ws.Range('E46').Value = '  -1.23%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '343.62'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -5.58%  '
$ws.Range('E48').Value = '  +3.25%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '32.76'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +6.51%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.01'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.31%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '6.34'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +3.78%  '
